# BRMSdata.xlsx edit: update BPR/BMR reference codes on the
# "BPRproductdetails" and "Productdetails" sheets, and move the active
# sheet/selection from BPRproductdetails to Productdetails.

$wb = $excel.ActiveWorkbook

# --- BPRproductdetails sheet -------------------------------------------------
$wsBPR = $wb.Worksheets.Item("BPRproductdetails")

# O2: BPR/3002022-04-00 -> " BPR/3002022-05-01"
$wsBPR.Range("O2").Value = " BPR/3002022-05-01"

# O3 is a new cell: " BPR/3002022-05-00"
$wsBPR.Range("O3").Value = " BPR/3002022-05-00"

# --- Productdetails sheet ----------------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdetails")

# C3: USA -> ROW
$wsProd.Range("C3").Value = "ROW"

# P2: BMR/2003021-09-01 -> BMR/2003021-13-00
$wsProd.Range("P2").Value = "BMR/2003021-13-00"

# O2: BMR/2003024-03-01 -> BMR/2003024-04-00
$wsProd.Range("O2").Value = "BMR/2003024-04-00"

# --- Active sheet / selection updates ---------------------------------------
# Previously BPRproductdetails was the active tab with E9 selected; now
# Productdetails is the active tab (with O2 selected), and
# BPRproductdetails's selection moves to H6.
[void]$wsBPR.Range("H6").Select()

[void]$wsProd.Activate()
[void]$wsProd.Range("O2").Select()
